$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2417.8572
$ws.Range("I70").Value = 1233.3334
$ws.Range("J70").Value = 3306.25
$ws.Range("K70").Value = 3700.0002
$ws.Range("L70").Value = 9918.75
$ws.Range("M70").Value = -3430.0002
$ws.Range("N70").Value = -10458.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 2417.8572
$ws.Range("I73").Value = 1233.3334
$ws.Range("J73").Value = 3306.25
$ws.Range("K73").Value = 3700.0002
$ws.Range("L73").Value = 9918.75
$ws.Range("M73").Value = -2764.0002
$ws.Range("N73").Value = -11790.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3171
$ws.Range("I132").Value = 3311.366
$ws.Range("K132").Value = 9934.098
$ws.Range("M132").Value = -7404.098

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 12928458
$ws.Range("I135").Value = 735.86957
$ws.Range("J135").Value = 29447214
$ws.Range("K135").Value = 6622.826129999999
$ws.Range("L135").Value = 265024926
$ws.Range("M135").Value = -4087.826129999999
$ws.Range("N135").Value = -265029996

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2135.516
$ws.Range("I138").Value = 1521.4822
$ws.Range("J138").Value = 3064.8647
$ws.Range("K138").Value = 4564.446599999999
$ws.Range("L138").Value = 9194.5941
$ws.Range("M138").Value = 575.5534000000007
$ws.Range("N138").Value = -19474.5941

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 964.76086
$ws.Range("I141").Value = 248.97144
$ws.Range("K141").Value = 746.91432
$ws.Range("M141").Value = 4433.08568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3390.33
$ws.Range("I32").Value = 3099.8538
$ws.Range("J32").Value = 4713.6113
$ws.Range("K32").Value = 3099.8538
$ws.Range("L32").Value = 4713.6113
$ws.Range("M32").Value = -2812.8538
$ws.Range("N32").Value = -5287.6113

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3004158
$ws.Range("I61").Value = 3087565
$ws.Range("J61").Value = 1514
$ws.Range("K61").Value = 3087565
$ws.Range("L61").Value = 1514
$ws.Range("M61").Value = -3087353
$ws.Range("N61").Value = -1938

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 9437441
$ws.Range("I74").Value = 12195837
$ws.Range("J74").Value = 12922.917
$ws.Range("K74").Value = 12195837
$ws.Range("L74").Value = 12922.917
$ws.Range("M74").Value = -12194963
$ws.Range("N74").Value = -14670.917

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 9437441
$ws.Range("I77").Value = 12195837
$ws.Range("J77").Value = 12922.917
$ws.Range("K77").Value = 60979185
$ws.Range("L77").Value = 64614.585
$ws.Range("M77").Value = -60974817
$ws.Range("N77").Value = -73350.58499999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 7034191
$ws.Range("I132").Value = 8359540.5
$ws.Range("J132").Value = 112922.22
$ws.Range("K132").Value = 25078621.5
$ws.Range("L132").Value = 338766.66
$ws.Range("M132").Value = -25076091.5
$ws.Range("N132").Value = -343826.66

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3004158
$ws.Range("I136").Value = 3087565
$ws.Range("J136").Value = 1514
$ws.Range("K136").Value = 9262695
$ws.Range("L136").Value = 4542
$ws.Range("M136").Value = -9260145
$ws.Range("N136").Value = -9642

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1394.5646
$ws.Range("I58").Value = 1056.4681
$ws.Range("J58").Value = 2453.9333
$ws.Range("K58").Value = 1056.4681
$ws.Range("L58").Value = 2453.9333
$ws.Range("M58").Value = -853.4681
$ws.Range("N58").Value = -2859.9333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1286.7297
$ws.Range("I134").Value = 1498.5769
$ws.Range("K134").Value = 4495.7307
$ws.Range("M134").Value = -1960.7307

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1394.5646
$ws.Range("I136").Value = 1056.4681
$ws.Range("J136").Value = 2453.9333
$ws.Range("K136").Value = 3169.4043
$ws.Range("L136").Value = 7361.7999
$ws.Range("M136").Value = -619.4043000000001
$ws.Range("N136").Value = -12461.7999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 7693042.5
$ws.Range("I5").Value = 347.96155
$ws.Range("J5").Value = 12821506
$ws.Range("K5").Value = 1043.88465
$ws.Range("L5").Value = 38464518
$ws.Range("M5").Value = -931.88465
$ws.Range("N5").Value = -38464742

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 10823778
$ws.Range("I122").Value = 19608488
$ws.Range("K122").Value = 176476392
$ws.Range("M122").Value = -176473942

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 7693042.5
$ws.Range("I135").Value = 347.96155
$ws.Range("J135").Value = 12821506
$ws.Range("K135").Value = 3131.65395
$ws.Range("L135").Value = 115393554
$ws.Range("M135").Value = -596.6539499999999
$ws.Range("N135").Value = -115398624

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 29413872
$ws.Range("I132").Value = 43480200
$ws.Range("J132").Value = 2464.7273
$ws.Range("K132").Value = 130440600
$ws.Range("L132").Value = 7394.1819
$ws.Range("M132").Value = -130438070
$ws.Range("N132").Value = -12454.1819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 191
$ws.Range("I55").Value = 112
$ws.Range("J55").Value = 234.88889
$ws.Range("K55").Value = 112
$ws.Range("L55").Value = 234.88889
$ws.Range("M55").Value = 61
$ws.Range("N55").Value = -580.8888899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2016.3
$ws.Range("I132").Value = 2014.5278
$ws.Range("J132").Value = 2032.25
$ws.Range("K132").Value = 6043.5834
$ws.Range("L132").Value = 6096.75
$ws.Range("M132").Value = -3513.5834
$ws.Range("N132").Value = -11156.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3000
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -4248

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 3000
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -21240

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6539300
$ws.Range("I132").Value = 7402607.5
$ws.Range("J132").Value = 2828.2856
$ws.Range("K132").Value = 22207822.5
$ws.Range("L132").Value = 8484.856800000001
$ws.Range("M132").Value = -22205292.5
$ws.Range("N132").Value = -13544.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3180180
$ws.Range("I136").Value = 9694.781999999999
$ws.Range("J136").Value = 6494778
$ws.Range("K136").Value = 29084.346
$ws.Range("L136").Value = 19484334
$ws.Range("M136").Value = -26534.346
$ws.Range("N136").Value = -19489434
